$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Rows 91 & 92 (1-based Excel rows) had their match data (columns
#    F:V) swapped - the "home"/"away" row order for the
#    sporting-san-jose/guanacasteca and cartagines/saprissa matches was
#    reversed. Columns A:E (index / pais / torneio / temporada / data)
#    stay untouched.
# ------------------------------------------------------------------

$row91 = @("Sporting San Jose", 1, "Guanacasteca", 1, 2.08, "19/10/2023 18:43", 2.05, "22/10/2023 00:40", 3.38, "19/10/2023 18:43", 3.47, "21/10/2023 19:13", 3.57, "19/10/2023 18:43", 3.7, "22/10/2023 00:40", "https://www.betexplorer.com/football/costa-rica/primera-division/sporting-san-jose-guanacasteca/UNdEriCI/")
$row92 = @("Cartagines", 0, "Saprissa", 4, 2.95, "19/10/2023 18:43", 3.05, "22/10/2023 00:50", 3.47, "19/10/2023 18:43", 3.65, "22/10/2023 00:50", 2.34, "19/10/2023 18:43", 2.26, "22/10/2023 00:50", "https://www.betexplorer.com/football/costa-rica/primera-division/cartagines-saprissa/vZ31ogSa/")

function Write-MatchCols {
    param($ws, $rowNum, $values)
    # $values has 17 entries mapping to columns F..V (6..22)
    for ($i = 0; $i -lt $values.Count; $i++) {
        $ws.Cells.Item($rowNum, 6 + $i).Value = $values[$i]
    }
}

# row 91 now gets what used to be row 92's data, and vice versa
Write-MatchCols $ws 91 $row92
Write-MatchCols $ws 92 $row91

# ------------------------------------------------------------------
# 2) Append four new match rows (128-131), mirroring the format of
#    the existing rows (row 127 is the current last data row).
# ------------------------------------------------------------------

$ws.Range("A127:V127").Copy()
$ws.Range("A128:V131").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$newRows = @(
    @{ R = 128; Idx = 127; E = 45256;      F = "Puntarenas FC"; G = 1; H = "Grecia";             I = 2; J = 1.66; K = "23/11/2023 22:12"; L = 1.75; M = "25/11/2023 23:56"; N = 3.92; O = "23/11/2023 22:12"; P = 3.79; Q = "25/11/2023 23:56"; R2 = 5;    S = "23/11/2023 22:12"; T = 4.68; U = "25/11/2023 23:56"; V = "https://www.betexplorer.com/football/costa-rica/primera-division/puntarenas-fc-grecia/phYuOJVR/" },
    @{ R = 129; Idx = 128; E = 45256.125;  F = "Alajuelense";  G = 3; H = "Guanacasteca";        I = 4; J = 1.39; K = "23/11/2023 02:12"; L = 1.5;  M = "26/11/2023 02:58"; N = 4.86; O = "23/11/2023 02:12"; P = 4.6;  Q = "26/11/2023 02:58"; R2 = 7.41; S = "23/11/2023 02:12"; T = 6.11; U = "26/11/2023 02:58"; V = "https://www.betexplorer.com/football/costa-rica/primera-division/alajuelense-guanacasteca/f3uTQLo9/" },
    @{ R = 130; Idx = 129; E = 45256.125;  F = "Herediano";    G = 3; H = "Sporting San Jose";  I = 0; J = 1.38; K = "23/11/2023 03:12"; L = 1.36; M = "26/11/2023 02:53"; N = 4.88; O = "23/11/2023 03:12"; P = 4.97; Q = "26/11/2023 02:53"; R2 = 7.61; S = "23/11/2023 03:12"; T = 8.539999999999999; U = "26/11/2023 02:53"; V = "https://www.betexplorer.com/football/costa-rica/primera-division/herediano-sporting-san-jose/rTyPR1W2/" },
    @{ R = 131; Idx = 130; E = 45256.125;  F = "Zeledon";      G = 1; H = "Cartagines";         I = 0; J = 3.13; K = "24/11/2023 02:12"; L = 3.55; M = "26/11/2023 02:51"; N = 3.59; O = "24/11/2023 02:12"; P = 3.71; Q = "26/11/2023 02:51"; R2 = 2.11; S = "24/11/2023 02:12"; T = 2.02; U = "26/11/2023 02:51"; V = "https://www.betexplorer.com/football/costa-rica/primera-division/zeledon-cartagines/GITyPaGL/" }
)

foreach ($row in $newRows) {
    $r = $row.R
    $ws.Cells.Item($r, 1).Value = $row.Idx
    $ws.Cells.Item($r, 2).Value = "costa-rica"
    $ws.Cells.Item($r, 3).Value = "primera-division"
    $ws.Cells.Item($r, 4).Value = "2023-2024"
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = $row.Q
    $ws.Cells.Item($r, 18).Value = $row.R2
    $ws.Cells.Item($r, 19).Value = $row.S
    $ws.Cells.Item($r, 20).Value = $row.T
    $ws.Cells.Item($r, 21).Value = $row.U
    $ws.Cells.Item($r, 22).Value = $row.V
}
